# CryCompanywiseStockReport_1.xlsx — correction pass.
#
# A batch of item rows had their Code/Rate/Value/Qty/Amount columns
# (B, D, E, F, G) shifted by one row within their item group during a
# previous import. This script rotates each affected group's data back:
# row i receives the data that belonged to row i+1 (wrapping around to
# row 0 for the last row in the group), while Sl.No (A) and the item
# Name (C) stay put (C is included in the rotation too, since it is
# effectively constant within a group except for one case-only spelling
# difference that also needs to move with its row).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$groups = @(
    @(149,150),
    @(161,162,163),
    @(264,265),
    @(279,280),
    @(313,314),
    @(316,317,318),
    @(346,347),
    @(350,351,352),
    @(372,373),
    @(375,376),
    @(379,380),
    @(382,383),
    @(389,390),
    @(419,420),
    @(421,422),
    @(431,432),
    @(536,537),
    @(579,580),
    @(599,600),
    @(601,602),
    @(687,688),
    @(859,860),
    @(889,890)
)

foreach ($g in $groups) {
    $n = $g.Length

    $bVals = @()
    $cVals = @()
    $dVals = @()
    $eVals = @()
    $fVals = @()
    $gVals = @()

    foreach ($r in $g) {
        $bVals += ,$ws.Cells.Item($r, 2).Value()
        $cVals += ,$ws.Cells.Item($r, 3).Value()
        $dVals += ,$ws.Cells.Item($r, 4).Value()
        $eVals += ,$ws.Cells.Item($r, 5).Value()
        $fVals += ,$ws.Cells.Item($r, 6).Value()
        $gVals += ,$ws.Cells.Item($r, 7).Value()
    }

    for ($i = 0; $i -lt $n; $i++) {
        $src = ($i + 1) % $n
        $r = $g[$i]
        $ws.Cells.Item($r, 2).Value = $bVals[$src]
        $ws.Cells.Item($r, 3).Value = $cVals[$src]
        $ws.Cells.Item($r, 4).Value = $dVals[$src]
        $ws.Cells.Item($r, 5).Value = $eVals[$src]
        $ws.Cells.Item($r, 6).Value = $fVals[$src]
        $ws.Cells.Item($r, 7).Value = $gVals[$src]
    }
}

Write-Host "Rotated $($groups.Length) item groups"
